$d = $word.ActiveDocument

# --- Change 1: rework the "Office of Community Control" monitoring sentence ---
# Old:  "Defendant shall report to the Office of Community Control forthwith for the following monitoring:"
# New:  "Prior to release the Defendant shall be fitted by the Office of Community Control for the following monitoring unit:"
$old1 = "Defendant shall report to the Office of Community Control forthwith for the following monitoring:"
$new1 = "Prior to release the Defendant shall be fitted by the Office of Community Control for the following monitoring unit:"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Change 2: refresh the "Vehicle Seizure/Immobilization" heading run ---
# (repagination moved the lastRenderedPageBreak hint off of this heading in the
# authored edit; re-touching the run brings the text back in sync)
$rngHeading = $d.Content
$found2 = $rngHeading.Find.Execute("Vehicle Seizure/Immobilization", $true, $false, $false, $false, $false, $true, 1, $false, "Vehicle Seizure/Immobilization", 2)

Write-Output ("change1=" + $found1 + " change2=" + $found2)
